$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (volume number, week-covering dates) ---
$ws.Range("A8").Value = "Volume 29   Number  48"
$ws.Range("C9").Value = "Report Covering the Week  11/28/2022  Through  12/4/2022"

# --- Weekly crime statistics table updates (rows 14-29) ---
# Use stable existing text cells as Copy() sources so that cells which must
# become the special text markers "0" (shared string) and "***.*" (shared
# string) pick up the exact same cell style/shared-string reference that
# Excel already uses elsewhere in the sheet, instead of Excel reinterpreting
# a typed "0" as a numeric value.
$zeroTextSrc = $ws.Range("C14")   # text "0"
$naTextSrc   = $ws.Range("E14")   # text "***.*"

$zeroTextSrc.Copy($ws.Range("G14"))
$naTextSrc.Copy($ws.Range("H14"))
$zeroTextSrc.Copy($ws.Range("C15"))
$zeroTextSrc.Copy($ws.Range("D15"))
$naTextSrc.Copy($ws.Range("E15"))
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 0
$ws.Range("N15").Value = -9.677419354838
$ws.Range("C16").Value = 8
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = 60
$ws.Range("F16").Value = 27
$ws.Range("G16").Value = 26
$ws.Range("H16").Value = 3.846153846153
$ws.Range("I16").Value = 288
$ws.Range("J16").Value = 183
$ws.Range("K16").Value = 57.377049180327
$ws.Range("L16").Value = 57.377049180327
$ws.Range("M16").Value = 7.063197026022
$ws.Range("N16").Value = -64.223602484472
$ws.Range("D17").Value = 9
$ws.Range("E17").Value = -33.333333333333
$ws.Range("F17").Value = 28
$ws.Range("G17").Value = 30
$ws.Range("H17").Value = -6.666666666666
$ws.Range("I17").Value = 399
$ws.Range("J17").Value = 295
$ws.Range("K17").Value = 35.254237288135
$ws.Range("L17").Value = 64.876033057851
$ws.Range("M17").Value = 9.315068493150
$ws.Range("N17").Value = -17.219917012448
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 18
$ws.Range("H18").Value = 63.636363636363
$ws.Range("I18").Value = 152
$ws.Range("J18").Value = 94
$ws.Range("K18").Value = 61.702127659574
$ws.Range("L18").Value = 8.571428571428
$ws.Range("M18").Value = -24
$ws.Range("N18").Value = -84.599797365754
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = -25
$ws.Range("F19").Value = 24
$ws.Range("G19").Value = 27
$ws.Range("H19").Value = -11.111111111111
$ws.Range("I19").Value = 370
$ws.Range("J19").Value = 305
$ws.Range("K19").Value = 21.311475409836
$ws.Range("L19").Value = 32.616487455197
$ws.Range("M19").Value = 46.825396825396
$ws.Range("N19").Value = -11.057692307692
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 29
$ws.Range("G20").Value = 13
$ws.Range("H20").Value = 123.076923076923
$ws.Range("I20").Value = 208
$ws.Range("J20").Value = 130
$ws.Range("K20").Value = 60
$ws.Range("L20").Value = 123.655913978495
$ws.Range("M20").Value = 188.888888888889
$ws.Range("N20").Value = -49.268292682926
$ws.Range("C21").Value = 28
$ws.Range("D21").Value = 30
$ws.Range("E21").Value = -6.666666666666
$ws.Range("F21").Value = 128
$ws.Range("G21").Value = 109
$ws.Range("H21").Value = 17.431192660550
$ws.Range("I21").Value = 1450
$ws.Range("J21").Value = 1047
$ws.Range("K21").Value = 38.490926456542
$ws.Range("L21").Value = 48.717948717948
$ws.Range("M21").Value = 22.673434856176
$ws.Range("N21").Value = -54.287515762925
$zeroTextSrc.Copy($ws.Range("D22"))
$naTextSrc.Copy($ws.Range("E22"))
$ws.Range("F22").Value = 4
$ws.Range("H22").Value = -60
$ws.Range("I22").Value = 29
$ws.Range("K22").Value = -3.333333333333
$ws.Range("M22").Value = 26.086956521739
$zeroTextSrc.Copy($ws.Range("F23"))
$ws.Range("G23").Value = 4
$ws.Range("H23").Value = -100
$ws.Range("J23").Value = 19
$ws.Range("K23").Value = -15.789473684210
$ws.Range("M23").Value = 33.333333333333
$ws.Range("C24").Value = 27
$ws.Range("D24").Value = 7
$ws.Range("E24").Value = 285.714285714286
$ws.Range("F24").Value = 85
$ws.Range("G24").Value = 59
$ws.Range("H24").Value = 44.067796610169
$ws.Range("I24").Value = 755
$ws.Range("J24").Value = 510
$ws.Range("K24").Value = 48.039215686274
$ws.Range("L24").Value = 19.841269841269
$ws.Range("M24").Value = 21.774193548387
$ws.Range("C25").Value = 12
$ws.Range("D25").Value = 11
$ws.Range("E25").Value = 9.090909090909
$ws.Range("F25").Value = 48
$ws.Range("G25").Value = 33
$ws.Range("H25").Value = 45.454545454545
$ws.Range("I25").Value = 511
$ws.Range("J25").Value = 369
$ws.Range("K25").Value = 38.482384823848
$ws.Range("L25").Value = 49.853372434017
$ws.Range("M25").Value = -17.313915857605
$zeroTextSrc.Copy($ws.Range("C26"))
$zeroTextSrc.Copy($ws.Range("D26"))
$naTextSrc.Copy($ws.Range("E26"))
$ws.Range("G26").Value = 3
$ws.Range("H26").Value = 33.333333333333
$ws.Range("C27").Value = 2
$zeroTextSrc.Copy($ws.Range("D27"))
$naTextSrc.Copy($ws.Range("E27"))
$ws.Range("F27").Value = 9
$ws.Range("G27").Value = 8
$ws.Range("H27").Value = 12.5
$ws.Range("I27").Value = 92
$ws.Range("K27").Value = -9.803921568627
$ws.Range("L27").Value = 70.370370370370
$zeroTextSrc.Copy($ws.Range("D28"))
$naTextSrc.Copy($ws.Range("E28"))
$ws.Range("G28").Value = 2
$ws.Range("L28").Value = 4.545454545454
$ws.Range("N28").Value = -79.824561403508
$zeroTextSrc.Copy($ws.Range("D29"))
$naTextSrc.Copy($ws.Range("E29"))
$ws.Range("G29").Value = 2
$ws.Range("L29").Value = -5
$ws.Range("N29").Value = -80
